$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Subject"
$ws.Range("B1").Value = "HbA1c"

$ws.Range("A2").Select()
